$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text values (e.g. "3.160.51" with dot-grouping)
# that look numeric to Excel auto-detection; force text format so they
# round-trip as strings, matching the source data.
$priceCells = @("D2","D3","D5","D6","D7","D8","D11","D15","D16","D17","D18","D20","D21","D23","D24","D25","D27","D28","D31","D33","D36","D38","D42","D43","D46","D47","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.075.60"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.160.51"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "574.69"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").Value = "150.04"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "3.158.87"
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").Value = "6.10"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("E13").Value = "  +12.54%  "
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").Value = "3.678.61"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "65.110.20"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "3.163.61"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").Value = "7.11"
$ws.Range("E18").Value = "  +2.17%  "
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").Value = "505.94"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").Value = "14.84"
$ws.Range("E21").Value = "  +1.90%  "
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("D23").Value = "15.29"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "7.73"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").Value = "84.38"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "8.91"
$ws.Range("E27").Value = "  +5.63%  "
$ws.Range("D28").Value = "2.91"
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("E30").Value = "  +5.46%  "
$ws.Range("D31").Value = "27.61"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "1.19"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("E34").Value = "  +3.94%  "
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("D36").Value = "54.91"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("E37").Value = "  +7.72%  "
$ws.Range("D38").Value = "467.52"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("E40").Value = "  +4.76%  "
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("D42").Value = "3.045.29"
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("D43").Value = "0.117"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("E44").Value = "  +5.61%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "28.49"
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("D47").Value = "0.0₃0589"
$ws.Range("E47").Value = "  +9.68%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("B50").Value = "CoreDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D50").Value = "2.57"
$ws.Range("E50").Value = "  +24.21%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "2.25"
$ws.Range("E51").Value = "  +2.26%  "
